$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add "yes" value to column C for the last three Yaaseen.Choudhury rows
$ws.Range("C14").Value = "yes"
$ws.Range("C15").Value = "yes"
$ws.Range("C16").Value = "yes"

# Update the active selection as recorded in the sheet view
$ws.Range("C13").Select()
